$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newText) {
    $c = $ws.Range($cellRef)
    $escaped = $newText.Replace('"', '""')
    $c.Formula = '="' + $escaped + '"'
    $c.Copy()
    $c.PasteSpecial(-4163)
}

Set-TextValue "H2" "8067.61"
Set-TextValue "H3" "17416.81"
Set-TextValue "H4" "9625.00"
Set-TextValue "H5" "249977.00"
Set-TextValue "H6" "5200.00"
Set-TextValue "H7" "105.28"
Set-TextValue "H8" "79870.00"
Set-TextValue "H9" "139782.56"
Set-TextValue "H10" "6088.92"
Set-TextValue "H11" "1085.50"
Set-TextValue "H12" "8134.12"
Set-TextValue "H13" "4530.00"
Set-TextValue "H14" "210.00"
Set-TextValue "H15" "713.70"
Set-TextValue "H16" "354.65"
Set-TextValue "H17" "469.29"
Set-TextValue "H18" "35295.00"
Set-TextValue "H19" "82.23"
Set-TextValue "H20" "24595.88"
Set-TextValue "H21" "164.00"
Set-TextValue "H22" "497.98"
Set-TextValue "H23" "146.06"
Set-TextValue "H24" "144.11"
Set-TextValue "H25" "22000.00"
Set-TextValue "H26" "18294.30"
Set-TextValue "H27" "1550.27"
Set-TextValue "H28" "50225.42"
Set-TextValue "H29" "39.99"
Set-TextValue "H30" "4199.18"
Set-TextValue "H31" "985.19"
Set-TextValue "H32" "2280.50"
Set-TextValue "H33" "310.67"
Set-TextValue "H34" "5164.10"
Set-TextValue "H35" "6716.20"
Set-TextValue "H36" "5132.39"
Set-TextValue "H37" "88.00"
Set-TextValue "H38" "100.00"
Set-TextValue "H39" "33739.38"
Set-TextValue "H40" "32.00"
Set-TextValue "H41" "3081.62"
Set-TextValue "H42" "5165.86"
Set-TextValue "H43" "520.39"
Set-TextValue "H44" "39000.00"
Set-TextValue "H45" "53460.00"
Set-TextValue "H46" "320.00"
Set-TextValue "H47" "3226.95"
Set-TextValue "H48" "43010.74"
Set-TextValue "H49" "243.16"
Set-TextValue "H50" "96.71"
Set-TextValue "H51" "5.64"
Set-TextValue "H52" "1140.00"
Set-TextValue "H53" "13390.00"
Set-TextValue "H54" "2631.36"
Set-TextValue "H55" "4866.00"
Set-TextValue "H56" "6875.00"
Set-TextValue "H57" "33.10"
Set-TextValue "H58" "414.52"
Set-TextValue "H59" "3585.00"
Set-TextValue "H60" "60.00"
Set-TextValue "H61" "182370.00"
Set-TextValue "H62" "2300.00"
Set-TextValue "H63" "3078.24"
Set-TextValue "H64" "60.92"
Set-TextValue "H65" "293.30"
Set-TextValue "H66" "4030.00"
Set-TextValue "H67" "719.00"
Set-TextValue "H68" "1449.00"
Set-TextValue "H69" "427.00"
Set-TextValue "H70" "157.32"
Set-TextValue "H71" "1377.15"
Set-TextValue "H72" "12.00"
Set-TextValue "H73" "150.20"
Set-TextValue "H74" "886.80"
Set-TextValue "H75" "37.35"
Set-TextValue "H76" "1843.00"
Set-TextValue "H77" "50.00"
Set-TextValue "H78" "450.00"
Set-TextValue "H79" "5040.00"
Set-TextValue "H80" "14216.00"
Set-TextValue "H81" "1398.57"
Set-TextValue "H82" "34.68"
Set-TextValue "H83" "1640.00"
Set-TextValue "H84" "2351.00"
Set-TextValue "H85" "11227.48"
Set-TextValue "H86" "8712.00"
Set-TextValue "H87" "500.00"
Set-TextValue "H88" "650.00"
Set-TextValue "H89" "100.00"
Set-TextValue "H90" "125.00"
Set-TextValue "H91" "1000.00"
Set-TextValue "H92" "100.00"
Set-TextValue "H93" "3500.00"
Set-TextValue "H94" "2031.16"
Set-TextValue "H95" "1368.00"
Set-TextValue "H96" "699.98"
Set-TextValue "H97" "1193.02"
Set-TextValue "H98" "17160.00"
Set-TextValue "H99" "154.00"
Set-TextValue "H100" "1491.07"
Set-TextValue "H101" "220.00"
Set-TextValue "H102" "10.36"
Set-TextValue "H103" "1846.96"
Set-TextValue "H104" "540.00"
Set-TextValue "H105" "350.00"
Set-TextValue "H106" "1694.00"
Set-TextValue "H107" "561.90"
Set-TextValue "H108" "0.81"
Set-TextValue "H109" "6044.09"
Set-TextValue "H110" "5823.00"
Set-TextValue "H111" "5625.14"
Set-TextValue "H112" "118.80"
Set-TextValue "H113" "38.00"
Set-TextValue "H114" "19259.00"
Set-TextValue "H115" "4272.00"
Set-TextValue "H116" "890.00"
Set-TextValue "H117" "1066.20"
Set-TextValue "H118" "1454.42"
Set-TextValue "H119" "2078.46"
Set-TextValue "H120" "993093.55"
Set-TextValue "H121" "980.00"
Set-TextValue "H122" "350.00"
Set-TextValue "H123" "30000.00"
Set-TextValue "H124" "17000.00"
Set-TextValue "H125" "50200.00"
Set-TextValue "H126" "10000.00"
Set-TextValue "H127" "33300.00"
Set-TextValue "H128" "32000.00"
Set-TextValue "H129" "5650.00"
Set-TextValue "H130" "6200.00"
Set-TextValue "H131" "278.00"
Set-TextValue "H132" "112.00"
Set-TextValue "H133" "4257.99"
Set-TextValue "E58" "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
Set-TextValue "F58" "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
Set-TextValue "E59" "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
Set-TextValue "E69" "TRABICHET MARIA. VERGARA ADEL Y OTRA"
Set-TextValue "F69" "TRABICHET MARIA. VERGARA ADEL Y OTRA"
